$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibition list)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 287
$ws1.Range("F3").Value = 1089
$ws1.Range("F4").Value = 2503
$ws1.Range("F5").Value = 212

# Sheet "全部类型" (all types combined list)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 287
$ws4.Range("F5").Value = 1089
$ws4.Range("F6").Value = 2503
$ws4.Range("F8").Value = 212
